$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture a plain ("General"-format, unstyled) cell style to restore after
# forcing Text format on numeric-looking strings, so only the displayed
# text changes and no cell keeps a lingering Text/quote-prefix style.
$plainStyle = $ws.Range("A1").Style

$ws.Range("D2").Value = "56.389.90"
$ws.Range("E2").Value = "  +3.92%  "

$ws.Range("D3").Value = "2.318.43"
$ws.Range("E3").Value = "  +2.38%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "519.46"
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = "  +4.75%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.13"
$ws.Range("D6").Style = $plainStyle
$ws.Range("E6").Value = "  +4.12%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.993"
$ws.Range("D7").Style = $plainStyle
$ws.Range("E7").Value = "  -0.44%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.537"
$ws.Range("D8").Style = $plainStyle
$ws.Range("E8").Value = "  +2.14%  "

$ws.Range("D9").Value = "2.342.38"
$ws.Range("E9").Value = "  +3.19%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.104"
$ws.Range("D10").Style = $plainStyle
$ws.Range("E10").Value = "  +8.62%  "

$ws.Range("E11").Value = "  +1.21%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.19"
$ws.Range("D12").Style = $plainStyle
$ws.Range("E12").Value = "  +7.64%  "

$ws.Range("E13").Value = "  +2.11%  "

$ws.Range("E14").Value = "  +5.04%  "

$ws.Range("D15").Value = "2.732.00"
$ws.Range("E15").Value = "  +2.52%  "

$ws.Range("D16").Value = "56.544.45"
$ws.Range("E16").Value = "  +4.26%  "

$ws.Range("E17").Value = "  +4.88%  "

$ws.Range("D18").Value = "2.327.51"
$ws.Range("E18").Value = "  +2.56%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.58"
$ws.Range("D19").Style = $plainStyle
$ws.Range("E19").Value = "  +3.39%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.30"
$ws.Range("D20").Style = $plainStyle
$ws.Range("E20").Value = "  +4.20%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "321.16"
$ws.Range("D21").Style = $plainStyle
$ws.Range("E21").Value = "  +5.93%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.65"
$ws.Range("D22").Style = $plainStyle
$ws.Range("E22").Value = "  +5.10%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = $plainStyle
$ws.Range("E23").Value = "  +0.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.90"
$ws.Range("D24").Style = $plainStyle
$ws.Range("E24").Value = "  +0.32%  "

$ws.Range("E25").Value = "  -0.47%  "

$ws.Range("E26").Value = "  +7.33%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.70"
$ws.Range("D27").Style = $plainStyle
$ws.Range("E27").Value = "  +5.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "171.44"
$ws.Range("D28").Style = $plainStyle
$ws.Range("E28").Value = "  +0.34%  "

$ws.Range("E29").Value = "  +12.42%  "

$ws.Range("D30").Value = "0.0₃0737"
$ws.Range("E30").Value = "  +6.81%  "

$ws.Range("E31").Value = "  +5.97%  "

$ws.Range("E32").Value = "  +5.17%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.45"
$ws.Range("D33").Style = $plainStyle
$ws.Range("E33").Value = "  +3.75%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("D34").Style = $plainStyle
$ws.Range("E34").Value = "  -0.04%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.990"
$ws.Range("D35").Style = $plainStyle
$ws.Range("E35").Value = "  -0.61%  "

$ws.Range("E36").Value = "  +5.92%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.943"
$ws.Range("D37").Style = $plainStyle
$ws.Range("E37").Value = "  +0.49%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.06"
$ws.Range("D38").Style = $plainStyle
$ws.Range("E38").Value = "  +9.40%  "

$ws.Range("E39").Value = "  +9.10%  "

$ws.Range("E40").Value = "  +4.48%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.382"
$ws.Range("D41").Style = $plainStyle
$ws.Range("E41").Value = "  +2.20%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "141.01"
$ws.Range("D42").Style = $plainStyle
$ws.Range("E42").Value = "  +13.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.61"
$ws.Range("D43").Style = $plainStyle
$ws.Range("E43").Value = "  +7.08%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "281.18"
$ws.Range("D44").Style = $plainStyle
$ws.Range("E44").Value = "  +16.37%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.16"
$ws.Range("D45").Style = $plainStyle
$ws.Range("E45").Value = "  +3.41%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0512"
$ws.Range("D46").Style = $plainStyle
$ws.Range("E46").Value = "  +3.78%  "

$ws.Range("E47").Value = "  +3.94%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.558"
$ws.Range("D48").Style = $plainStyle
$ws.Range("E48").Value = "  +2.46%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0217"
$ws.Range("D49").Style = $plainStyle
$ws.Range("E49").Value = "  +6.14%  "

$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.38"
$ws.Range("D51").Style = $plainStyle
$ws.Range("E51").Value = "  +7.80%  "
